$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.685507
$ws.Range("H2").Value = 11.056521
$ws.Range("I2").Value = 0.3585631737883472
$ws.Range("J2").Value = 0.3585631737883472
$ws.Range("O2").Value = 0.921725411846598
$ws.Range("P2").Value = 0.9217254118465981
$ws.Range("Q2").Value = 23.38916231227567
$ws.Range("R2").Value = 210.502460810481
$ws.Range("S2").Value = 0.3304967890330877
$ws.Range("T2").Value = 0.3304967890330877

# Row 3
$ws.Range("G3").Value = 3.685507
$ws.Range("H3").Value = 11.056521
$ws.Range("I3").Value = 0.3585631737883472
$ws.Range("J3").Value = 0.3585631737883472
$ws.Range("M3").Value = 0.5389353333333333
$ws.Range("N3").Value = 1.616806
$ws.Range("O3").Value = 0.07827458815340194
$ws.Range("P3").Value = 0.07827458815340194
$ws.Range("Q3").Value = 1.986249943547333
$ws.Range("R3").Value = 17.876249491926
$ws.Range("S3").Value = 0.02806638475525956
$ws.Range("T3").Value = 0.02806638475525956

# Row 4
$ws.Range("I4").Value = 0.009647184430711629
$ws.Range("J4").Value = 0.009647184430711629
$ws.Range("O4").Value = 0.921725411846598
$ws.Range("P4").Value = 0.9217254118465981
$ws.Range("S4").Value = 0.008892055042557764
$ws.Range("T4").Value = 0.008892055042557766

# Row 5
$ws.Range("I5").Value = 0.009647184430711629
$ws.Range("J5").Value = 0.009647184430711629
$ws.Range("M5").Value = 0.5389353333333333
$ws.Range("N5").Value = 1.616806
$ws.Range("O5").Value = 0.07827458815340194
$ws.Range("P5").Value = 0.07827458815340194
$ws.Range("Q5").Value = 0.053440288718
$ws.Range("R5").Value = 0.480962598462
$ws.Range("S5").Value = 0.0007551293881538641
$ws.Range("T5").Value = 0.0007551293881538641

# Row 6
$ws.Range("G6").Value = 6.493877
$ws.Range("H6").Value = 19.481631
$ws.Range("I6").Value = 0.6317896417809412
$ws.Range("J6").Value = 0.6317896417809411
$ws.Range("O6").Value = 0.921725411846598
$ws.Range("P6").Value = 0.9217254118465981
$ws.Range("Q6").Value = 41.21179072213233
$ws.Range("R6").Value = 370.906116499191
$ws.Range("S6").Value = 0.5823365677709527
$ws.Range("T6").Value = 0.5823365677709526

# Row 7
$ws.Range("G7").Value = 6.493877
$ws.Range("H7").Value = 19.481631
$ws.Range("I7").Value = 0.6317896417809412
$ws.Range("J7").Value = 0.6317896417809411
$ws.Range("M7").Value = 0.5389353333333333
$ws.Range("N7").Value = 1.616806
$ws.Range("O7").Value = 0.07827458815340194
$ws.Range("P7").Value = 0.07827458815340194
$ws.Range("Q7").Value = 3.499779765620667
$ws.Range("R7").Value = 31.498017890586
$ws.Range("S7").Value = 0.04945307400998852
$ws.Range("T7").Value = 0.0494530740099885
